$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B20").Value = 6185
$ws.Range("D20").Value = 5578287
$ws.Range("E20").Value = 901.9057396928051
$ws.Range("F20").Value = 6.840559682155822
$ws.Range("H20").Value = 26.18637478953996
